$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "66.997.42"
$ws.Range("E2").Value = "  +1.22%  "

# Row 3
$ws.Range("D3").Value = "3.119.09"
$ws.Range("E3").Value = "  +1.37%  "

# Row 4
$ws.Range("E4").Value = "  -0.01%  "

# Row 5
$ws.Range("D5").Value = "576.80"
$ws.Range("E5").Value = "  -0.31%  "

# Row 6
$ws.Range("D6").Value = "173.54"
$ws.Range("E6").Value = "  +3.60%  "

# Row 7
$ws.Range("E7").Value = "  -0.02%  "

# Row 8
$ws.Range("D8").Value = "3.114.26"
$ws.Range("E8").Value = "  +1.28%  "

# Row 9
$ws.Range("D9").Value = "0.522"
$ws.Range("E9").Value = "  -0.09%  "

# Row 10
$ws.Range("D10").Value = "6.44"
$ws.Range("E10").Value = "  -3.50%  "

# Row 11
$ws.Range("D11").Value = "0.154"
$ws.Range("E11").Value = "  +0.67%  "

# Row 12
$ws.Range("D12").Value = "0.478"
$ws.Range("E12").Value = "  -1.34%  "

# Row 13
$ws.Range("D13").Value = "0.0000247"
$ws.Range("E13").Value = "  -0.65%  "

# Row 14
$ws.Range("D14").Value = "37.16"
$ws.Range("E14").Value = "  +1.02%  "

# Row 15
$ws.Range("E15").Value = "  -0.98%  "

# Row 16
$ws.Range("D16").Value = "3.637.03"
$ws.Range("E16").Value = "  +1.38%  "

# Row 17
$ws.Range("D17").Value = "67.004.85"
$ws.Range("E17").Value = "  +1.03%  "

# Row 18
$ws.Range("D18").Value = "7.12"
$ws.Range("E18").Value = "  -1.21%  "

# Row 19
$ws.Range("D19").Value = "3.117.05"
$ws.Range("E19").Value = "  +1.17%  "

# Row 20
$ws.Range("D20").Value = "16.24"
$ws.Range("E20").Value = "  -0.90%  "

# Row 21
$ws.Range("D21").Value = "477.90"
$ws.Range("E21").Value = "  +2.59%  "

# Row 22
$ws.Range("D22").Value = "0.712"
$ws.Range("E22").Value = "  -0.10%  "

# Row 23
$ws.Range("E23").Value = "  +3.03%  "

# Row 24
$ws.Range("D24").Value = "83.89"
$ws.Range("E24").Value = "  +0.73%  "

# Row 25
$ws.Range("D25").Value = "13.24"
$ws.Range("E25").Value = "  +2.77%  "

# Row 26
$ws.Range("D26").Value = "2.30"
$ws.Range("E26").Value = "  +1.27%  "

# Row 27
$ws.Range("B27").Value = "Dai"
$ws.Range("C27").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  -0.03%  "

# Row 28
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").Value = "9.98"
$ws.Range("E28").Value = "  -0.55%  "

# Row 29
$ws.Range("D29").Value = "7.98"
$ws.Range("E29").Value = "  -1.83%  "

# Row 30
$ws.Range("D30").Value = "2.41"
$ws.Range("E30").Value = "  -0.19%  "

# Row 31
$ws.Range("E31").Value = "  +0.07%  "

# Row 32
$ws.Range("D32").Value = "28.70"
$ws.Range("E32").Value = "  +1.41%  "

# Row 33
$ws.Range("D33").Value = "0.0₃0973"
$ws.Range("E33").Value = "  -3.87%  "

# Row 34
$ws.Range("D34").Value = "0.113"
$ws.Range("E34").Value = "  -2.99%  "

# Row 35
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  -0.20%  "

# Row 36
$ws.Range("D36").Value = "5.86"
$ws.Range("E36").Value = "  -0.32%  "

# Row 37
$ws.Range("D37").Value = "0.982"
$ws.Range("E37").Value = "  -1.33%  "

# Row 38
$ws.Range("D38").Value = "47.73"
$ws.Range("E38").Value = "  -0.22%  "

# Row 39
$ws.Range("D39").Value = "2.08"
$ws.Range("E39").Value = "  +2.82%  "

# Row 40
$ws.Range("D40").Value = "50.03"
$ws.Range("E40").Value = "  +0.39%  "

# Row 41
$ws.Range("D41").Value = "0.310"
$ws.Range("E41").Value = "  -1.25%  "

# Row 42
$ws.Range("D42").Value = "0.122"
$ws.Range("E42").Value = "  +0.86%  "

# Row 43
$ws.Range("D43").Value = "8.61"
$ws.Range("E43").Value = "  -0.68%  "

# Row 44
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "2.812.03"
$ws.Range("E44").Value = "  +1.86%  "

# Row 45
$ws.Range("B45").Value = "dogwifhat"
$ws.Range("C45").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D45").Value = "2.60"
$ws.Range("E45").Value = "  -8.79%  "

# Row 46
$ws.Range("D46").Value = "0.0355"
$ws.Range("E46").Value = "  -1.66%  "

# Row 47
$ws.Range("D47").Value = "379.92"
$ws.Range("E47").Value = "  -0.58%  "

# Row 48
$ws.Range("D48").Value = "135.85"
$ws.Range("E48").Value = "  +1.21%  "

# Row 49
$ws.Range("E49").Value = "  +0.00%  "

# Row 50
$ws.Range("D50").Value = "24.73"
$ws.Range("E50").Value = "  +1.19%  "

# Row 51
$ws.Range("E51").Value = "  -0.72%  "
